$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = '2019 Maruti Swift LXI'
$ws.Range("C1").Value = '31,466 KM'
$ws.Range("D1").Value = 'MANUAL'
$ws.Range("E1").Value = 'REGULARLY SERVICED'
$ws.Range("F1").Value = '₹5.14 Lakh'
$ws.Range("G1").Value = 'EMIs from ₹10,044/month'

$ws.Range("B2").Value = '2017 Maruti Wagon R 1.0 LXI'
$ws.Range("C2").Value = '79,416 KM'
$ws.Range("D2").Value = 'MANUAL'
$ws.Range("E2").Value = 'LOW RUNNING COST'
$ws.Range("F2").Value = '₹3.30 Lakh'
$ws.Range("G2").Value = 'EMIs from ₹6,452/month'

$ws.Range("B3").Value = '2022 Maruti Vitara Brezza ZXI PLUS'
$ws.Range("C3").Value = '20,059 KM'
$ws.Range("D3").Value = 'MANUAL'
$ws.Range("E3").Value = 'TOP MODEL'
$ws.Range("F3").Value = '₹8.70 Lakh'
$ws.Range("G3").Value = 'EMIs from ₹16,560/month'

$ws.Range("B4").Value = '2018 Maruti Alto 800 LXI'
$ws.Range("C4").Value = '89,267 KM'
$ws.Range("D4").Value = 'MANUAL'
$ws.Range("E4").Value = 'LOW RUNNING COST'
$ws.Range("F4").Value = '₹2.78 Lakh'
$ws.Range("G4").Value = 'EMIs from ₹5,435/month'

$ws.Range("B5").Value = '2019 Maruti Swift ZXI PLUS AMT'
$ws.Range("C5").Value = '36,023 KM'
$ws.Range("D5").Value = 'AUTOMATIC'
$ws.Range("E5").Value = 'TOP MODEL'
$ws.Range("F5").Value = '₹6.31 Lakh'
$ws.Range("G5").Value = 'EMIs from ₹12,336/month'

$ws.Range("B6").Value = '2018 Maruti Celerio VXI'
$ws.Range("C6").Value = '52,238 KM'
$ws.Range("D6").Value = 'MANUAL'
$ws.Range("E6").Value = '100% TYRE LIFE REMAINING'
$ws.Range("F6").Value = '₹3.98 Lakh'
$ws.Range("G6").Value = 'EMIs from ₹7,781/month'

$ws.Range("B7").Value = '2012 Maruti Wagon R 1.0 VXI'
$ws.Range("C7").Value = '30,121 KM'
$ws.Range("D7").Value = 'MANUAL'
$ws.Range("E7").Value = 'TOP MODEL'
$ws.Range("F7").Value = '₹2.07 Lakh'
$ws.Range("G7").Value = 'EMIs from ₹6,875/month'

$ws.Range("B8").Value = '2022 Maruti Celerio VXI CNG'
$ws.Range("C8").Value = '6,790 KM'
$ws.Range("D8").Value = 'MANUAL'
$ws.Range("E8").Value = 'REGULARLY SERVICED'
$ws.Range("F8").Value = '₹6.20 Lakh'
$ws.Range("G8").Value = 'EMIs from ₹12,121/month'

$ws.Range("B9").Value = '2012 Maruti Swift Dzire VXI'
$ws.Range("C9").Value = '22,466 KM'
$ws.Range("D9").Value = 'MANUAL'
$ws.Range("E9").Value = 'REGULARLY SERVICED'
$ws.Range("F9").Value = '₹2.91 Lakh'
$ws.Range("G9").Value = 'EMIs from ₹9,665/month'

$ws.Range("B10").Value = '2016 Maruti Baleno ZETA PETROL 1.2'
$ws.Range("C10").Value = '96,466 KM'
$ws.Range("D10").Value = 'MANUAL'
$ws.Range("E10").Value = 'ALLOY WHEELS'
$ws.Range("F10").Value = '₹4.59 Lakh'
$ws.Range("G10").Value = 'EMIs from ₹8,974/month'

$ws.Range("B11").Value = '2022 Maruti Vitara Brezza VXI AT SHVS'
$ws.Range("C11").Value = '1,402 KM'
$ws.Range("D11").Value = 'AUTOMATIC'
$ws.Range("E11").Value = 'STANDARD SAFETY FEATURES'
$ws.Range("F11").Value = '₹10.64 Lakh'
$ws.Range("G11").Value = 'EMIs from ₹20,252/month'

$ws.Range("B12").Value = '2018 Maruti Dzire VXI AMT'
$ws.Range("C12").Value = '45,657 KM'
$ws.Range("D12").Value = 'AUTOMATIC'
$ws.Range("E12").Value = 'REGULARLY SERVICED'
$ws.Range("F12").Value = '₹5.24 Lakh'
$ws.Range("G12").Value = 'EMIs from ₹10,244/month'

$ws.Range("B13").Value = '2017 Maruti Dzire LXI'
$ws.Range("C13").Value = '22,148 KM'
$ws.Range("D13").Value = 'MANUAL'
$ws.Range("E13").Value = 'STANDARD SAFETY FEATURES'
$ws.Range("F13").Value = '₹5.13 Lakh'
$ws.Range("G13").Value = 'EMIs from ₹10,029/month'

$ws.Range("B14").Value = '2017 Maruti Swift ZXI'
$ws.Range("C14").Value = '63,367 KM'
$ws.Range("D14").Value = 'MANUAL'
$ws.Range("E14").Value = 'TOP MODEL'
$ws.Range("F14").Value = '₹5.43 Lakh'
$ws.Range("G14").Value = 'EMIs from ₹10,616/month'

$ws.Range("B15").Value = '2022 Maruti Swift ZXI PLUS'
$ws.Range("C15").Value = '24,626 KM'
$ws.Range("D15").Value = 'MANUAL'
$ws.Range("E15").Value = 'TOP MODEL'
$ws.Range("F15").Value = '₹8.01 Lakh'
$ws.Range("G15").Value = 'EMIs from ₹15,246/month'

$ws.Range("B16").Value = '2010 Maruti Wagon R 1.0 VXI'
$ws.Range("C16").Value = '84,625 KM'
$ws.Range("D16").Value = 'MANUAL'
$ws.Range("E16").Value = 'TOP MODEL'
$ws.Range("F16").Value = '₹1.42 Lakh'
$ws.Range("G16").Value = 'EMIs from ₹12,617/month'

$ws.Range("B17").Value = '2022 Maruti Baleno ZETA PETROL 1.2'
$ws.Range("C17").Value = '15,515 KM'
$ws.Range("D17").Value = 'MANUAL'
$ws.Range("E17").Value = 'ALLOY WHEELS'
$ws.Range("F17").Value = '₹8.61 Lakh'
$ws.Range("G17").Value = 'EMIs from ₹16,388/month'

$ws.Range("B18").Value = '2022 Maruti S PRESSO VXI (O) CNG'
$ws.Range("C18").Value = '21,240 KM'
$ws.Range("D18").Value = 'MANUAL'
$ws.Range("E18").Value = 'REGULARLY SERVICED'
$ws.Range("F18").Value = '₹4.75 Lakh'
$ws.Range("G18").Value = 'EMIs from ₹9,286/month'

$ws.Range("B19").Value = '2019 Maruti Baleno ZETA PETROL 1.2'
$ws.Range("C19").Value = '30,890 KM'
$ws.Range("D19").Value = 'MANUAL'
$ws.Range("E19").Value = 'ALLOY WHEELS'
$ws.Range("F19").Value = '₹5.89 Lakh'
$ws.Range("G19").Value = 'EMIs from ₹11,515/month'

$ws.Range("B20").Value = '2013 Maruti Wagon R 1.0 VXI'
$ws.Range("C20").Value = '48,692 KM'
$ws.Range("D20").Value = 'MANUAL'
$ws.Range("E20").Value = 'TOP MODEL'
$ws.Range("F20").Value = '₹2.79 Lakh'
$ws.Range("G20").Value = 'EMIs from ₹7,347/month'
